$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text (inline string) nature instead of being
# auto-converted to a number by Excel when values look numeric (e.g. "14.90").
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '38.345.75'
$ws.Range("E2").Value = '  +3.75%  '

# Row 3
$ws.Range("D3").Value = '2.068.60'
$ws.Range("E3").Value = '  +3.02%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").Value = '230.88'
$ws.Range("E5").Value = '  +2.19%  '

# Row 6
$ws.Range("D6").Value = '0.618'
$ws.Range("E6").Value = '  +2.25%  '

# Row 7
$ws.Range("D7").Value = '61.54'
$ws.Range("E7").Value = '  +11.66%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '0.389'
$ws.Range("E9").Value = '  +4.78%  '

# Row 10
$ws.Range("E10").Value = '  +4.88%  '

# Row 11
$ws.Range("E11").Value = '  +2.34%  '

# Row 12
$ws.Range("D12").Value = '14.90'
$ws.Range("E12").Value = '  +6.35%  '

# Row 13
$ws.Range("D13").Value = '2.372.73'
$ws.Range("E13").Value = '  +2.95%  '

# Row 14
$ws.Range("D14").Value = '21.53'
$ws.Range("E14").Value = '  +9.33%  '

# Row 15
$ws.Range("D15").Value = '0.767'
$ws.Range("E15").Value = '  +4.70%  '

# Row 16
$ws.Range("D16").Value = '5.35'
$ws.Range("E16").Value = '  +4.31%  '

# Row 17
$ws.Range("D17").Value = '2.070.44'
$ws.Range("E17").Value = '  +3.31%  '

# Row 18
$ws.Range("D18").Value = '38.244.38'
$ws.Range("E18").Value = '  +3.65%  '

# Row 19
$ws.Range("E19").Value = '  +2.75%  '

# Row 20
$ws.Range("D20").Value = '70.26'
$ws.Range("E20").Value = '  +3.01%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0839'
$ws.Range("E21").Value = '  +3.68%  '

# Row 22
$ws.Range("D22").Value = '226.34'
$ws.Range("E22").Value = '  +1.70%  '

# Row 23
$ws.Range("E23").Value = '  +0.00%  '

# Row 24
$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  -0.48%  '

# Row 25
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  +4.63%  '

# Row 26
$ws.Range("D26").Value = '9.36'
$ws.Range("E26").Value = '  +5.02%  '

# Row 27
$ws.Range("D27").Value = '166.53'
$ws.Range("E27").Value = '  +1.28%  '

# Row 28
$ws.Range("D28").Value = '0.133'
$ws.Range("E28").Value = '  +6.19%  '

# Row 29
$ws.Range("D29").Value = '19.17'
$ws.Range("E29").Value = '  +3.45%  '

# Row 30
$ws.Range("E30").Value = '  +2.46%  '

# Row 31
$ws.Range("E31").Value = '  +3.05%  '

# Row 32
$ws.Range("D32").Value = '4.58'
$ws.Range("E32").Value = '  +4.52%  '

# Row 33
$ws.Range("D33").Value = '4.66'
$ws.Range("E33").Value = '  +5.19%  '

# Row 34
$ws.Range("E34").Value = '  +10.06%  '

# Row 35
$ws.Range("D35").Value = '0.0609'
$ws.Range("E35").Value = '  +1.67%  '

# Row 36
$ws.Range("D36").Value = '2.33'
$ws.Range("E36").Value = '  +0.62%  '

# Row 37
$ws.Range("D37").Value = '6.26'
$ws.Range("E37").Value = '  +17.30%  '

# Row 38
$ws.Range("D38").Value = '3.35'
$ws.Range("E38").Value = '  +6.33%  '

# Row 39
$ws.Range("E39").Value = '  +0.10%  '

# Row 40
$ws.Range("D40").Value = '1.528.46'
$ws.Range("E40").Value = '  +4.52%  '

# Row 41
$ws.Range("D41").Value = '17.26'
$ws.Range("E41").Value = '  +9.15%  '

# Row 42
$ws.Range("D42").Value = '98.49'
$ws.Range("E42").Value = '  +4.59%  '

# Row 43
$ws.Range("D43").Value = '0.0219'
$ws.Range("E43").Value = '  +3.88%  '

# Row 44
$ws.Range("E44").Value = '  +4.21%  '

# Row 45
$ws.Range("D45").Value = '0.0934'
$ws.Range("E45").Value = '  +2.79%  '

# Row 46
$ws.Range("E46").Value = '  +1.90%  '

# Row 47
$ws.Range("E47").Value = '  -2.88%  '

# Row 48
$ws.Range("E48").Value = '  +3.41%  '

# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '7.11'
$ws.Range("E49").Value = '  +1.00%  '

# Row 50
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '2.94'
$ws.Range("E50").Value = '  +2.06%  '

# Row 51
$ws.Range("D51").Value = '2.260.18'
$ws.Range("E51").Value = '  +2.97%  '
